# Docx writer: Use different style for block quotes in notes.
#
# Pandoc's docx writer used to render a block quote inside a footnote
# with the regular "Block Text" style. That style is based on "Body
# Text", which doesn't make sense inside a note (wrong font size, etc).
# Give block quotes in notes their own "Footnote Block Text" style
# instead: same paragraph formatting as "Block Text" (spacing before/
# after + left/right indent), but based on "Footnote Text" so it picks
# up the footnote font size/next-style instead of the body one.
#
# Closes #9243.

$d = $word.ActiveDocument

$styleName = "Footnote Block Text"

$alreadyExists = $false
try {
    $null = $d.Styles.Item($styleName)
    $alreadyExists = $true
} catch {
    $alreadyExists = $false
}

if (-not $alreadyExists) {
    $blockText = $d.Styles.Add($styleName, $wdStyleTypeParagraph)

    # Same base/next chain as "Block Text", but rooted in "Footnote
    # Text" rather than "Body Text".
    $blockText.BaseStyle = "Footnote Text"
    $blockText.NextParagraphStyle = "Footnote Text"

    $blockText.Priority = 9
    $blockText.UnhideWhenUsed = $true
    $blockText.QuickStyle = $true

    # Mirror "Block Text"'s paragraph formatting: 100 twips (5pt)
    # spacing before/after, 480 twips (24pt) left/right indent, no
    # first-line indent.
    $blockText.ParagraphFormat.SpaceBefore = 5
    $blockText.ParagraphFormat.SpaceAfter = 5
    $blockText.ParagraphFormat.FirstLineIndent = 0
    $blockText.ParagraphFormat.LeftIndent = 24
    $blockText.ParagraphFormat.RightIndent = 24
}
